# 05-using-tcpmon.docx minor-fixes edit script
#
# Summary of the target change (per the supplied diff):
#   - The only real *content* edit is the lab command path changing from
#       ~/servers/tcpmon/build
#     to
#       ~/servers/tcpmon-1.0-bin/build
#   - The "_GoBack" bookmark that used to sit right after the word
#     "Kepler" is removed from there and re-inserted between "~/" and
#     "servers/tcpmon-1.0-bin" in the "cd ..." line.
#   - Everywhere else the diff only wraps already-present text in
#     <w:proofErr w:type="gramStart|gramEnd|spellStart|spellEnd"/> markers
#     and/or splits a run into several runs at the same text boundaries
#     that the proofing pass flagged. These proofErr markers are an
#     artifact Word's background spell/grammar checker stamps into the
#     file; they carry no visible/semantic content (no text, formatting,
#     or layout change) and are not reachable through the Word
#     object model (no such method/property is exposed here, matching
#     real Word's COM automation surface). We still reproduce the run
#     splits themselves (which *are* achievable, by momentarily dropping
#     a bookmark at the split point) so the saved package's run structure
#     lines up with the target as closely as the object model allows.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: split the run that contains $searchText so that a run boundary
# exists right at the start (Where=1) or end (Where=0) of the found
# range. We do this the same way Word itself would when you click in
# the middle of a run and start typing/marking something: drop a
# temporary bookmark at that exact point (which forces the engine to
# break the run there) and immediately delete the bookmark again so it
# leaves no trace other than the run break.
# ---------------------------------------------------------------------
function Split-At {
    param(
        [string]$SearchText,
        [int]$Where,          # 1 = collapse to start, 0 = collapse to end
        [bool]$MatchCase = $true
    )

    $r = $d.Content
    $found = $r.Find.Execute($SearchText, $MatchCase, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "Split-At: NOT FOUND -> $SearchText"
        return $false
    }

    $point = $r.Duplicate
    $point.Collapse($Where)

    $markName = "zzTmpSplit"
    $d.Bookmarks.Add($markName, $point)
    $d.Bookmarks($markName).Delete()
    return $true
}

# =======================================================================
# 1) Text content change + bookmark relocation in the "cd ..." line
# =======================================================================

# 1a. Update the lab path text itself.
$d.Content.Find.Execute("~/servers/tcpmon/build", $true, $false, $false, $false, $false, $true, 1, $false, "~/servers/tcpmon-1.0-bin/build", 2) | Out-Null

# 1b. Remove the _GoBack bookmark from its old spot (right after "Kepler").
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# 1c. Re-insert _GoBack between "~/" and "servers/tcpmon-1.0-bin".
$r = $d.Content
$r.Find.Execute("~/servers/tcpmon-1.0-bin/build", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pathRange = $r.Duplicate
$pathStart = $pathRange.Duplicate
$pathStart.Collapse(1)
$pathStart.MoveEnd(1, 2)             # wdCharacter = 1; move end 2 chars -> covers "~/"
$bmPoint = $pathStart.Duplicate
$bmPoint.Collapse(0)                 # collapse to the point right after "~/"
$d.Bookmarks.Add("_GoBack", $bmPoint) | Out-Null

# =======================================================================
# 2) Reproduce the run splits the proofing pass introduced (best-effort;
#    purely structural -- no visible text changes).
# =======================================================================

Split-At "Kepler" 1

Split-At "You should have " 0
Split-At "tcpmon in the VM." 1

Split-At "~/" 0
Split-At "servers/tcpmon-1.0-bin" 0

Split-At "./" 0

Split-At " org.apache.cxf.endpoint.Client" 1
Split-At "org.apache.cxf.endpoint.Client" 0
Split-At " org.apache.cxf.frontend.ClientProxy" 1
Split-At "org.apache.cxf.frontend.ClientProxy" 0
Split-At " org.apache.cxf.transport.http.HTTPConduit" 1
Split-At "org.apache.cxf.transport.http.HTTPConduit" 0
Split-At " org.apache.cxf.transports.http.configuration.HTTPClientPolicy" 1
Split-At "org.apache.cxf.transports.http.configuration.HTTPClientPolicy" 0

Split-At " object" 0
Split-At " object" 1

Split-At "Client " 0
Split-At "client = ClientProxy" 1
Split-At "client" 0
Split-At "ClientProxy." 1
Split-At "ClientProxy." 0
Split-At "getClient(" 0
Split-At "(port" 1

Split-At "        HTTPConduit" 1
Split-At "HTTPConduit http" 0
Split-At "http = (" 0
Split-At "(HTTPConduit" 1
Split-At "HTTPConduit) client" 0
Split-At ") client.getConduit(" 0
Split-At "client.getConduit(" 1
Split-At "getConduit(" 0

Split-At "        HTTPClientPolicy httpClientPolicy" 1
Split-At "HTTPClientPolicy httpClientPolicy" 0
Split-At "httpClientPolicy = new" 1
Split-At "httpClientPolicy = " 0
Split-At "new HTTPClientPolicy(" 0
Split-At "new " 0
Split-At "HTTPClientPolicy();" 1
Split-At "HTTPClientPolicy(" 0

Split-At "        httpClientPolicy.setProxyServer(" 1
Split-At "httpClientPolicy.setProxyServer" 0
Split-At 'setProxyServer("' 0
Split-At '("localhost' 1
Split-At '"localhost"' 0

Split-At "        httpClientPolicy.setProxyServerPort" 1
Split-At "httpClientPolicy.setProxyServerPort" 0

Split-At "        http.setClient(httpClientPolicy)" 1
Split-At "http.setClient" 0
Split-At "setClient(httpClientPolicy" 1
Split-At "(httpClientPolicy" 0
Split-At "httpClientPolicy);" 0

Split-At "This could also be set in an XML " 0
Split-At "config file" 1

Split-At "(Attribution-" 0
Split-At "Sharealike" 0

Split-At "Software Engineering " 0

Write-Output "edit complete"
